$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.854.82"
$ws.Range("E2").Value = "  +0.68%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.699.29"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.73"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("E6").Value = "  -0.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4033"
$ws.Range("E7").Value = "  +2.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4060"
$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.006"
$ws.Range("E9").Value = "  -0.32%  "

$ws.Range("E10").Value = "  +2.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.462"
$ws.Range("E11").Value = "  -3.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08819"
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.69"
$ws.Range("E13").Value = "  +3.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.471"
$ws.Range("E14").Value = "  -1.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.024"
$ws.Range("E15").Value = "  -0.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001347"
$ws.Range("E16").Value = "  -1.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.723.33"
$ws.Range("E17").Value = "  +1.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.48"
$ws.Range("E18").Value = "  -2.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07176"
$ws.Range("E19").Value = "  +0.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.92"
$ws.Range("E20").Value = "  +5.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.213"
$ws.Range("E21").Value = "  -2.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  -0.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.57"
$ws.Range("E23").Value = "  +1.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.854.12"
$ws.Range("E24").Value = "  +0.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.337"
$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.877"
$ws.Range("E26").Value = "  -6.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.617"
$ws.Range("E27").Value = "  +26.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.03"
$ws.Range("E28").Value = "  +1.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.49"
$ws.Range("E29").Value = "  -0.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "143.84"
$ws.Range("E30").Value = "  +4.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.131"
$ws.Range("E31").Value = "  -4.47%  "

$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.977.86"
$ws.Range("E32").Value = "  +5.19%  "

$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.280"
$ws.Range("E33").Value = "  +13.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08738"
$ws.Range("E34").Value = "  -1.53%  "

$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.375"
$ws.Range("E35").Value = "  -2.08%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03188"
$ws.Range("E36").Value = "  +9.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.030"
$ws.Range("E37").Value = "  -1.65%  "

$ws.Range("E38").Value = "  +4.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8499"
$ws.Range("E39").Value = "  +8.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.89"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09418"
$ws.Range("E41").Value = "  +2.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.00"
$ws.Range("E42").Value = "  -2.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.468"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.76"
$ws.Range("E44").Value = "  +6.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.727"
$ws.Range("E45").Value = "  +5.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7449"
$ws.Range("E46").Value = "  +3.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.245"
$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.405"
$ws.Range("E48").Value = "  +5.87%  "

$ws.Range("E49").Value = "  -0.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.38"
$ws.Range("E50").Value = "  +1.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08359"
$ws.Range("E51").Value = "  +4.59%  "
